$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update currency data from Dolar/USD/$ to Iene/JPY/¥
$ws.Range("A2").Value = "Iene"
$ws.Range("B2").Value = "JPY"
$ws.Range("C2").Value = [char]0x00A5
$ws.Range("D2").Value = 0.0498
$ws.Range("E2").Value = 0.0492
$ws.Range("F2").Value = "Diminuiu 1.20%"

# Row 5: updated report time
$ws.Range("F5").Value = "21:59"
